$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Additional Metadata" column (F)
# so that column shifts right to G, and the newly inserted column F becomes
# the "participants" column (sds 3.0 performances.xlsx add participants column).
[void]$ws.Columns("F").Insert()

# Header text for the new column.
$ws.Range("F1").Value = "participants"

# Give the new column a sensible width, close to the original template's.
$ws.Columns("F").ColumnWidth = 10

# Reflect the active-cell/selection state recorded after the edit.
[void]$ws.Range("F2").Select()
